$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.378.18'
$ws.Range("E2").Value = '  +9.91%  '
$ws.Range("D3").Value = '1.675.80'
$ws.Range("E4").Value = '  -0.40%  '
$s = $ws.Range("D5").Style
$ws.Range("D5").Value = "'0.9990"
$ws.Range("D5").Style = $s
$ws.Range("E5").Value = '  +0.38%  '
$s = $ws.Range("D6").Style
$ws.Range("D6").Value = "'305.22"
$ws.Range("D6").Style = $s
$ws.Range("E6").Value = '  +2.51%  '
$s = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.3687"
$ws.Range("D7").Style = $s
$ws.Range("E7").Value = '  +1.90%  '
$s = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.3422"
$ws.Range("D8").Style = $s
$ws.Range("E8").Value = '  +2.67%  '
$s = $ws.Range("D9").Style
$ws.Range("D9").Value = "'47.56"
$ws.Range("D9").Style = $s
$ws.Range("E9").Value = '  +15.31%  '
$s = $ws.Range("D10").Style
$ws.Range("D10").Value = "'1.158"
$ws.Range("D10").Style = $s
$ws.Range("E10").Value = '  +4.43%  '
$s = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.07203"
$ws.Range("D11").Style = $s
$ws.Range("E11").Value = '  +3.95%  '
$ws.Range("E12").Value = '  -0.64%  '
$s = $ws.Range("D13").Style
$ws.Range("D13").Value = "'6.125"
$ws.Range("D13").Style = $s
$ws.Range("E13").Value = '  +5.71%  '
$s = $ws.Range("D14").Style
$ws.Range("D14").Value = "'20.10"
$ws.Range("D14").Style = $s
$ws.Range("E14").Value = '  +4.46%  '
$s = $ws.Range("D15").Style
$ws.Range("D15").Value = "'6.723"
$ws.Range("D15").Style = $s
$ws.Range("E15").Value = '  +3.21%  '
$ws.Range("D16").Value = '1.675.32'
$ws.Range("E16").Value = '  +5.69%  '
$s = $ws.Range("D17").Style
$ws.Range("D17").Value = "'0.00001101"
$ws.Range("D17").Style = $s
$ws.Range("E17").Value = '  +4.14%  '
$s = $ws.Range("D18").Style
$ws.Range("D18").Value = "'0.9991"
$ws.Range("D18").Style = $s
$ws.Range("E18").Value = '  +0.25%  '
$s = $ws.Range("D19").Style
$ws.Range("D19").Value = "'0.06644"
$ws.Range("D19").Style = $s
$ws.Range("E19").Value = '  +1.01%  '
$ws.Range("E20").Value = '  +5.76%  '
$s = $ws.Range("D21").Style
$ws.Range("D21").Value = "'16.45"
$ws.Range("D21").Style = $s
$ws.Range("E21").Value = '  +4.55%  '
$s = $ws.Range("D22").Style
$ws.Range("D22").Value = "'6.091"
$ws.Range("D22").Style = $s
$ws.Range("E22").Value = '  +3.35%  '
$s = $ws.Range("D23").Style
$ws.Range("D23").Value = "'12.16"
$ws.Range("D23").Style = $s
$ws.Range("E23").Value = '  +4.98%  '
$ws.Range("D24").Value = '24.312.47'
$ws.Range("E24").Value = '  +9.57%  '
$s = $ws.Range("D25").Style
$ws.Range("D25").Value = "'2.440"
$ws.Range("D25").Style = $s
$ws.Range("E25").Value = '  +2.74%  '
$s = $ws.Range("D26").Style
$ws.Range("D26").Value = "'2.651"
$ws.Range("D26").Style = $s
$ws.Range("E26").Value = '  +7.52%  '
$s = $ws.Range("D27").Style
$ws.Range("D27").Value = "'152.37"
$ws.Range("D27").Style = $s
$ws.Range("E27").Value = '  +3.08%  '
$s = $ws.Range("D28").Style
$ws.Range("D28").Value = "'19.37"
$ws.Range("D28").Style = $s
$ws.Range("E28").Value = '  +1.42%  '
$ws.Range("D29").Value = '1.860.49'
$ws.Range("E29").Value = '  +6.09%  '
$s = $ws.Range("D30").Style
$ws.Range("D30").Value = "'127.46"
$ws.Range("D30").Style = $s
$ws.Range("E30").Value = '  +5.57%  '
$s = $ws.Range("D31").Style
$ws.Range("D31").Value = "'6.270"
$ws.Range("D31").Style = $s
$ws.Range("E31").Value = '  +7.76%  '
$s = $ws.Range("D32").Style
$ws.Range("D32").Value = "'4.048"
$ws.Range("D32").Style = $s
$ws.Range("E32").Value = '  +2.68%  '
$s = $ws.Range("D33").Style
$ws.Range("D33").Value = "'0.9678"
$ws.Range("D33").Style = $s
$ws.Range("E33").Value = '  +6.19%  '
$s = $ws.Range("D34").Style
$ws.Range("D34").Value = "'0.08457"
$ws.Range("D34").Style = $s
$ws.Range("E34").Value = '  +4.24%  '
$s = $ws.Range("D35").Style
$ws.Range("D35").Value = "'1.674"
$ws.Range("D35").Style = $s
$ws.Range("E35").Value = '  +2.67%  '
$s = $ws.Range("D36").Style
$ws.Range("D36").Value = "'12.35"
$ws.Range("D36").Style = $s
$ws.Range("E36").Value = '  +6.14%  '
$s = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.06380"
$ws.Range("D37").Style = $s
$ws.Range("E37").Value = '  +7.28%  '
$s = $ws.Range("D38").Style
$ws.Range("D38").Value = "'5.298"
$ws.Range("D38").Style = $s
$ws.Range("E38").Value = '  +4.31%  '
$ws.Range("E39").Value = '  +6.47%  '
$s = $ws.Range("D40").Style
$ws.Range("D40").Value = "'8.660"
$ws.Range("D40").Style = $s
$ws.Range("E40").Value = '  +4.53%  '
$s = $ws.Range("D41").Style
$ws.Range("D41").Value = "'1.227"
$ws.Range("D41").Style = $s
$s = $ws.Range("D42").Style
$ws.Range("D42").Value = "'0.2083"
$ws.Range("D42").Style = $s
$ws.Range("E42").Value = '  +5.57%  '
$s = $ws.Range("D43").Style
$ws.Range("D43").Value = "'0.6071"
$ws.Range("D43").Style = $s
$ws.Range("E43").Value = '  +5.48%  '
$s = $ws.Range("D44").Style
$ws.Range("D44").Value = "'0.9988"
$ws.Range("D44").Style = $s
$ws.Range("E44").Value = '  +0.29%  '
$s = $ws.Range("D45").Style
$ws.Range("D45").Value = "'3.749"
$ws.Range("D45").Style = $s
$ws.Range("E45").Value = '  -0.37%  '
$ws.Range("E46").Value = '  +0.92%  '
$s = $ws.Range("D47").Style
$ws.Range("D47").Value = "'0.5860"
$ws.Range("D47").Style = $s
$ws.Range("E47").Value = '  +5.95%  '
$s = $ws.Range("D48").Style
$ws.Range("D48").Value = "'125.53"
$ws.Range("D48").Style = $s
$ws.Range("E48").Value = '  +0.74%  '
$s = $ws.Range("D49").Style
$ws.Range("D49").Value = "'2.012"
$ws.Range("D49").Style = $s
$ws.Range("E49").Value = '  +4.40%  '
$s = $ws.Range("D50").Style
$ws.Range("D50").Value = "'0.07136"
$ws.Range("D50").Style = $s
$ws.Range("E50").Value = '  +6.19%  '
$s = $ws.Range("D51").Style
$ws.Range("D51").Value = "'75.61"
$ws.Range("D51").Style = $s
$ws.Range("E51").Value = '  +4.96%  '
